$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row (A1:Q3 -> A1:Q4)
$null = $lo.ListRows.Add()

function Set-CellLikeRow3 {
    param($colLetter, $value)
    $ws.Range("$colLetter`3").Copy()
    $ws.Range("$colLetter`4").PasteSpecial(-4122)
    $ws.Range("$colLetter`4").Value = $value
}

# A4: timestamp
Set-CellLikeRow3 "A" 46013.42023517361

# B4: submitter email
Set-CellLikeRow3 "B" "mariobohorquezg@gmail.com"

# C4: Nombre del Apostolado
Set-CellLikeRow3 "C" "Movímiento de Espiritualidad Matrimonial (MEM)"

# D4: Instagram
Set-CellLikeRow3 "D" "memmaracaibo "

# E4: Templo (reuses existing shared string "San Ramón Nonato")
Set-CellLikeRow3 "E" "San Ramón Nonato"

# F4: Descripcion del Apostolado
Set-CellLikeRow3 "F" "Fomentar la espiritualidad conyugal de los esposos y la santificación matrimonial a través de las actividades
ordinarias. Ser un factor de gran influencia en nuestras familias y en la sociedad para impulsar la vocación
matrimonial como célula de nuestra Iglesia.
Proveer un espacio para que los matrimonios revivan el amor conyugal, teniendo a Cristo como centro de su
alianza matrimonial a través de Retiros de Espiritualidad Matrimonial, la meditación de La Palabra, actividades
de formación y apoyando las actividades de apostolado de la parroquia."

# G4: Descripcion de actividades
Set-CellLikeRow3 "G" "Este movimiento nació en octubre de 1990 cuando el Padre Antonio Abella, entonces párroco de la iglesia San
Ramón Nonato, por inspiración divina conformó una serie de meditaciones sobre la vocación del matrimonio a
la luz de la Palabra de Dios y realizó el primer Retiro de Espiritualidad Matrimonial (REM) a tres matrimonios de
la parroquia.
El retiro se diseñó con la finalidad de intensificar la doctrina matrimonial y fundamentar la espiritualidad de los
esposos en la Palabra de Dios, mediante un conjunto de charlas, reflexiones, oraciones y celebraciones litúrgicas
a las que se recurre durante cuarenta y ocho horas que dura el retiro.
A partir del segundo retiro realizado en Maracaibo, el REM se constituyó como grupo parroquial, estableciéndose
reuniones semanales como medio fundamental para acrecentar la fe, canalizar la formación y el crecimiento
espiritual de los esposos. Posteriormente y junto a la realización de retiros periódicos y la participación de las
parejas en actividades apostólicas de la parroquia, el grupo se dedicó a dictar los cursillos prematrimoniales y
catequesis de niños y jóvenes.
Luego de consolidarse el REM como grupo parroquial en la década de los años 90, por iniciativa del Padre Antonio
Abella y con la aceptación y beneplácito de todos los matrimonios participantes, en octubre del 2000 el grupo
pasó a llamarse Movimiento de Espiritualidad Matrimonial (MEM).
La sucesión de la dirección espiritual del MEM en sentido cronológico desde su nacimiento, es como sigue:
• P. Antonio Abella (1990-1996)
• P. José Zaporta (1996-1998)
• P. Francisco Ortiz (1998-2000)
• P. Jesús García (2000-2003)
• P. Néstor Burgos (2007-2009)
• P. Jesús Bel (2009-2014)
• P. Richard Godoy (2014 hasta el presente)
Con el pasar del tiempo las actividades del REM fueron paulatinamente pasando al liderazgo de los propios
matrimonios del movimiento quienes se han pasado el testigo manteniendo intactas las reflexiones originales
del P. Antonio Abella, siempre con la guía espiritual del párroco de San Ramón Nonato y bajo el amparo y gracia
del Espíritu Santo a quien consagramos con fervor nuestro movimiento y apostolado.
Desde 2018 el MEM se constituyó como la sección “Adultos y Familia” de la Acción Católica de Venezuela y
comenzó a realizar retiros de espiritualidad matrimonial dedicados a nuevas parroquias de Maracaibo,
empezando por Santa Rosa de Lima y San Bartolomé Apóstol (Ziruma), para convertirse en una institución laical
enfocada a la Pastoral Familiar, que tiene como principal apostolado la lucha por el fortalecimiento del
matrimonio católico."

# H4: Nombre del Coordinador 1
Set-CellLikeRow3 "H" "Mario Bohórquez de Fernández "

# I4: Telefono del Coordinador 1 (stored as a plain number, like the column's other cells)
$ws.Range("B3").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 4246144593

# J4: Nombre del Coordinador 2
Set-CellLikeRow3 "J" "Mónica Fernández de Bohorquez "

# K4: Telefono del Coordinador 2 (kept as text, quote-prefixed like K3)
Set-CellLikeRow3 "K" "'04146178604"

# L4: Quienes pueden participar
Set-CellLikeRow3 "L" "Matrimonios o parejas de unión libre en concubinato o parejas de divorciados y vueltos a casar "

# M4: Requisitos (reuses existing shared string "No")
Set-CellLikeRow3 "M" "No"

# N4: dias que se reunen
Set-CellLikeRow3 "N" "Miércoles"

# O4: horario
Set-CellLikeRow3 "O" "6:30 pm - 8:30 pm "

# P4: Foto de portada (hyperlink)
Set-CellLikeRow3 "P" "https://drive.google.com/open?id=1EZ875Kc0bR40h4zJ1FqrcCJuoj9LQqYv"
$ws.Hyperlinks.Add($ws.Range("P4"), "https://drive.google.com/open?id=1EZ875Kc0bR40h4zJ1FqrcCJuoj9LQqYv")
$ws.Range("P3").Copy()
$ws.Range("P4").PasteSpecial(-4122)

# Q4: Fotos adicionales (hyperlink)
Set-CellLikeRow3 "Q" "https://drive.google.com/open?id=1Ag50xlJ6e4U5FaG71pOVC_KifmN2od49"
$ws.Hyperlinks.Add($ws.Range("Q4"), "https://drive.google.com/open?id=1Ag50xlJ6e4U5FaG71pOVC_KifmN2od49")
$ws.Range("Q3").Copy()
$ws.Range("Q4").PasteSpecial(-4122)

# Match the row height used by the other data rows (set last: wrapped long
# text in F4/G4 would otherwise auto-grow the row after this point).
$ws.Rows.Item(4).RowHeight = 22.5

Write-Host "Row 4 populated; table range:" $lo.Range.Address()
